# Update TPM-derived statistics for Sema4d-Erbb2 LR pairs (YoungD0)
# Reflects rerun of NATMI scripts with new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1349983333333333
$ws.Range("H2").Value = 0.404995
$ws.Range("I2").Value = 0.06188478316908706
$ws.Range("J2").Value = 0.06188478316908706
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.572065666666667
$ws.Range("N2").Value = 4.716197
$ws.Range("O2").Value = 0.1759712293834306
$ws.Range("P2").Value = 0.1759712293834305
$ws.Range("Q2").Value = 0.2122262448905556
$ws.Range("R2").Value = 1.910036204015
$ws.Range("S2").Value = 0.01088994137439128
$ws.Range("T2").Value = 0.01088994137439128

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1349983333333333
$ws.Range("H3").Value = 0.404995
$ws.Range("I3").Value = 0.06188478316908706
$ws.Range("J3").Value = 0.06188478316908706
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.237801
$ws.Range("N3").Value = 12.713403
$ws.Range("O3").Value = 0.4743638053196239
$ws.Range("P3").Value = 0.4743638053196239
$ws.Range("Q3").Value = 0.5720960719983333
$ws.Range("R3").Value = 5.148864647985
$ws.Range("S3").Value = 0.02935590123546795
$ws.Range("T3").Value = 0.02935590123546795

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1349983333333333
$ws.Range("H4").Value = 0.404995
$ws.Range("I4").Value = 0.06188478316908706
$ws.Range("J4").Value = 0.06188478316908706
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.123785
$ws.Range("N4").Value = 9.371355
$ws.Range("O4").Value = 0.3496649652969456
$ws.Range("P4").Value = 0.3496649652969455
$ws.Range("Q4").Value = 0.4217057686916666
$ws.Range("R4").Value = 3.795351918225
$ws.Range("S4").Value = 0.02163894055922783
$ws.Range("T4").Value = 0.02163894055922782

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.9347409999999999
$ws.Range("H5").Value = 2.804223
$ws.Range("I5").Value = 0.4284959871424753
$ws.Range("J5").Value = 0.4284959871424753
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.572065666666667
$ws.Range("N5").Value = 4.716197
$ws.Range("O5").Value = 0.1759712293834306
$ws.Range("P5").Value = 0.1759712293834305
$ws.Range("Q5").Value = 1.469474233325667
$ws.Range("R5").Value = 13.225268099931
$ws.Range("S5").Value = 0.07540296564332805
$ws.Range("T5").Value = 0.07540296564332803

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.9347409999999999
$ws.Range("H6").Value = 2.804223
$ws.Range("I6").Value = 0.4284959871424753
$ws.Range("J6").Value = 0.4284959871424753
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.237801
$ws.Range("N6").Value = 12.713403
$ws.Range("O6").Value = 0.4743638053196239
$ws.Range("P6").Value = 0.4743638053196239
$ws.Range("Q6").Value = 3.961246344541
$ws.Range("R6").Value = 35.651217100869
$ws.Range("S6").Value = 0.2032629870250932
$ws.Range("T6").Value = 0.2032629870250932

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.9347409999999999
$ws.Range("H7").Value = 2.804223
$ws.Range("I7").Value = 0.4284959871424753
$ws.Range("J7").Value = 0.4284959871424753
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.123785
$ws.Range("N7").Value = 9.371355
$ws.Range("O7").Value = 0.3496649652969456
$ws.Range("P7").Value = 0.3496649652969455
$ws.Range("Q7").Value = 2.919929914685
$ws.Range("R7").Value = 26.279369232165
$ws.Range("S7").Value = 0.1498300344740541
$ws.Range("T7").Value = 0.149830034474054

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.111707
$ws.Range("H8").Value = 3.335121
$ws.Range("I8").Value = 0.5096192296884376
$ws.Range("J8").Value = 0.5096192296884376
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.572065666666667
$ws.Range("N8").Value = 4.716197
$ws.Range("O8").Value = 0.1759712293834306
$ws.Range("P8").Value = 0.1759712293834305
$ws.Range("Q8").Value = 1.747676406093
$ws.Range("R8").Value = 15.729087654837
$ws.Range("S8").Value = 0.08967832236571126
$ws.Range("T8").Value = 0.08967832236571124

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.111707
$ws.Range("H9").Value = 3.335121
$ws.Range("I9").Value = 0.5096192296884376
$ws.Range("J9").Value = 0.5096192296884376
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.237801
$ws.Range("N9").Value = 12.713403
$ws.Range("O9").Value = 0.4743638053196239
$ws.Range("P9").Value = 0.4743638053196239
$ws.Range("Q9").Value = 4.711193036307
$ws.Range("R9").Value = 42.40073732676301
$ws.Range("S9").Value = 0.2417449170590627
$ws.Range("T9").Value = 0.2417449170590627

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.111707
$ws.Range("H10").Value = 3.335121
$ws.Range("I10").Value = 0.5096192296884376
$ws.Range("J10").Value = 0.5096192296884376
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.123785
$ws.Range("N10").Value = 9.371355
$ws.Range("O10").Value = 0.3496649652969456
$ws.Range("P10").Value = 0.3496649652969455
$ws.Range("Q10").Value = 3.472733650995
$ws.Range("R10").Value = 31.254602858955
$ws.Range("S10").Value = 0.1781959902636637
$ws.Range("T10").Value = 0.1781959902636636
